$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.877.13"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.083.76"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.66"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.37"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.27"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "2.027.84"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "37.777.36"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.54"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").Value = "  +4.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.40"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.03"
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.137"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.74"
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0986"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.34"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.66"
$ws.Range("E43").Value = "  +6.78%  "
$ws.Range("D44").Value = "1.446.64"
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.20"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.41"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "2.275.11"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.88"
$ws.Range("E51").Value = "  +1.16%  "
